$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 10: E10 gets "COMPLETED" (keeps its existing style s="2")
$ws.Range("E10").Value = "COMPLETED"

# Row 11: restyle A11:D11 to match row 10's formatting (s=2 for text cols,
# s=3 for the date col) by copying row 10's formats over, then add the new
# E11 = "COMPLETED" cell with that same style.
$ws.Range("A10:D10").Copy()
$ws.Range("A11:D11").PasteSpecial(-4122)

$ws.Range("E11").Value = "COMPLETED"
$ws.Range("E10").Copy()
$ws.Range("E11").PasteSpecial(-4122)

$excel.CutCopyMode = 0
